$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Unprotect()

# Update the confidential disclaimer text date: 2021-05-14 -> 2021-05-17
$ws.Range("A7").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + [char]10 + "Model holdings provided as of 2021-05-17 for illustrative purposes only and are subject to change."

# Update weight/percent-change figures on rows 2-4
$ws.Range("D2").Value = 0.8492741246797609
$ws.Range("E2").Value = -0.001885369532428571

$ws.Range("D3").Value = 0.1507258753202391
$ws.Range("E3").Value = -0.0009442870632673239

$ws.Range("E4").Value = -0.001743524053515744

$ws.Protect()
